$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.902.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.87%  "

$ws.Range("E9").Value = "  +3.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.90%  "

$ws.Range("E12").Value = "  +2.48%  "

$ws.Range("E13").Value = "  +2.63%  "

$ws.Range("E14").Value = "  +3.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.624.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("E16").Value = "  +3.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.293.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.825.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.26%  "

$ws.Range("E20").Value = "  +9.10%  "

$ws.Range("E21").Value = "  +2.36%  "

$ws.Range("E22").Value = "  +2.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("E25").Value = "  +4.07%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +5.08%  "

$ws.Range("E28").Value = "  +4.38%  "

$ws.Range("E29").Value = "  +3.43%  "

$ws.Range("E30").Value = "  -3.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.40%  "

$ws.Range("E35").Value = "  +5.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.39%  "

$ws.Range("E39").Value = "  +5.32%  "

$ws.Range("E40").Value = "  +2.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.68%  "

$ws.Range("E42").Value = "  +5.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.074.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("E46").Value = "  +3.55%  "

$ws.Range("E47").Value = "  +7.19%  "

$ws.Range("E48").Value = "  +4.34%  "

$ws.Range("E49").Value = "  +3.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.02%  "

$ws.Range("E51").Value = "  +3.56%  "
